$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.729.12"
$ws.Range("E2").Value = "  +0.18%  "

# Row 3
$ws.Range("D3").Value = "1.640.00"
$ws.Range("E3").Value = "  -0.19%  "

# Row 4
$ws.Range("E4").Value = "  +0.11%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.14"
$ws.Range("E5").Value = "  +1.00%  "

# Row 6
$ws.Range("E6").Value = "  -0.27%  "

# Row 7
$ws.Range("E7").Value = "  +0.19%  "

# Row 8
$ws.Range("E8").Value = "  -0.14%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0626"
$ws.Range("E9").Value = "  -0.05%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.13"
$ws.Range("E10").Value = "  +0.36%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0842"
$ws.Range("E11").Value = "  -0.22%  "

# Row 12
$ws.Range("D12").Value = "1.869.76"
$ws.Range("E12").Value = "  -0.11%  "

# Row 13
$ws.Range("D13").Value = "1.642.96"
$ws.Range("E13").Value = "  +0.31%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.15"
$ws.Range("E14").Value = "  -0.58%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.526"
$ws.Range("E15").Value = "  -0.47%  "

# Row 16
$ws.Range("E16").Value = "  -0.95%  "

# Row 17
$ws.Range("D17").Value = "26.729.06"
$ws.Range("E17").Value = "  +0.11%  "

# Row 18
$ws.Range("D18").Value = "0.0₃0733"
$ws.Range("E18").Value = "  -1.35%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "213.92"
$ws.Range("E19").Value = "  -1.13%  "

# Row 20
$ws.Range("E20").Value = "  +0.14%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.38"
$ws.Range("E22").Value = "  +5.69%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.21"
$ws.Range("E23").Value = "  -0.68%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.26"
$ws.Range("E24").Value = "  -2.54%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.44"
$ws.Range("E25").Value = "  -0.01%  "

# Row 26
$ws.Range("E26").Value = "  +0.33%  "

# Row 27
$ws.Range("E27").Value = "  -1.27%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.17"
$ws.Range("E28").Value = "  +0.06%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.59"
$ws.Range("E29").Value = "  -0.70%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0508"
$ws.Range("E30").Value = "  -1.42%  "

# Row 31
$ws.Range("E31").Value = "  +1.12%  "

# Row 32
$ws.Range("E32").Value = "  +0.25%  "

# Row 33
$ws.Range("E33").Value = "  -1.13%  "

# Row 34
$ws.Range("D34").Value = "1.284.51"
$ws.Range("E34").Value = "  +0.44%  "

# Row 35
$ws.Range("E35").Value = "  -0.27%  "

# Row 36
$ws.Range("E36").Value = "  +1.27%  "

# Row 37
$ws.Range("E37").Value = "  -0.82%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.535"
$ws.Range("E38").Value = "  +0.56%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.816"
$ws.Range("E39").Value = "  -1.43%  "

# Row 40
$ws.Range("E40").Value = "  +0.11%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.804"
$ws.Range("E41").Value = "  -1.06%  "

# Row 42
$ws.Range("E42").Value = "  -1.29%  "

# Row 43
$ws.Range("E43").Value = "  -2.74%  "

# Row 44
$ws.Range("D44").Value = "1.780.46"
$ws.Range("E44").Value = "  -0.08%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.26"
$ws.Range("E45").Value = "  +3.47%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.89"
$ws.Range("E46").Value = "  +0.07%  "

# Row 47
$ws.Range("E47").Value = "  +0.26%  "

# Row 48
$ws.Range("E48").Value = "  +0.31%  "

# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.65"
$ws.Range("E49").Value = "  -1.44%  "

# Row 50
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0966"
$ws.Range("E50").Value = "  +0.22%  "

# Row 51
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.406"
$ws.Range("E51").Value = "  -0.06%  "
